# UPDATE technology portfolios for Norway
$wb = $excel.ActiveWorkbook

# Update the base values on the "2025" sheet; the other sheets (2030-2050)
# reference '2025'!B2 via formulas and will recalculate automatically.
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 324328
$ws2025.Range("C2").Value = 11.4383

# The C2 "cop" value on the other sheets is a hard-coded (non-formula) value,
# so it must be updated explicitly on each sheet as well.
$sheetNames = @("2030", "2035", "2040", "2045", "2050")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C2").Value = 11.4383
}

$excel.CalculateFullRebuild()
